# Update countries & provincias Spain
# - Re-sort Armenia above Malasia/Marruecos (Armenia's case count overtook them)
# - Refresh the case-count snapshot for rows 60-62
# - Bump the "datos actualizados" timestamp from 08:35 to 09:05

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header (row 1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 09:05"

# Row 60: now Armenia, with updated totals
$ws.Range("A60").Value = "Armenia"
$ws.Range("B60").Value = 7774
$ws.Range("C60").Value = 372
$ws.Range("D60").Value = 3255
$ws.Range("E60").Value = 4421
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 7
$ws.Range("H60").Value = 98

# Row 61: now Malasia (previous Armenia-row pushed down one), unchanged totals
$ws.Range("A61").Value = "Malasia"
$ws.Range("B61").Value = 7604
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 6041
$ws.Range("E61").Value = 1448
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 115

# Row 62: now Marruecos, unchanged totals
$ws.Range("A62").Value = "Marruecos"
$ws.Range("B62").Value = 7577
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 4881
$ws.Range("E62").Value = 2494
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 202
